# SUGAR Unity Localization workbook update:
# Add a new "REMEMBER" / "Remember Me" localization entry, inserted as a new
# row just above the "SEARCH" row (i.e. it becomes the new row 35), pushing
# all subsequent rows (SEARCH ... ALLIANCES) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 35. This shifts the existing rows
# 35-43 (SEARCH ... ALLIANCES) down to rows 36-44.
$ws.Rows.Item(35).Insert() | Out-Null

# A freshly inserted row doesn't carry the sheet's standard custom row
# height, so restore it to match every other data row.
$ws.Rows.Item(35).RowHeight = 15.75

# The new entry is formatted like the other "header-ish" rows in this sheet
# (e.g. row 1, and the GROUP_MEMBERS / ALLIANCES rows), so copy that cell
# formatting into the new row's key/value cells.
$ws.Range("A1:B1").Copy() | Out-Null
$ws.Range("A35:B35").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new localization key/value pair.
$ws.Range("A35").Value = "REMEMBER"
$ws.Range("B35").Value = "Remember Me"

# Reflect the new active selection (the cursor now rests on B35 instead of
# the old B34) as last left by the editor.
$ws.Range("B35").Select() | Out-Null
